# Apply cryptos list update (price + volume refresh, row 50/51 swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.748.88"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "3.382.10"
$ws.Range("E3").Value = "  -3.84%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.89"
$ws.Range("E5").Value = "  -3.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.31"
$ws.Range("E6").Value = "  -6.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "3.381.23"
$ws.Range("E8").Value = "  -3.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.119"
$ws.Range("E10").Value = "  -10.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.02"
$ws.Range("E11").Value = "  -10.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.369"
$ws.Range("E12").Value = "  -8.12%  "

$ws.Range("D13").Value = "3.969.18"
$ws.Range("E13").Value = "  -3.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000175"
$ws.Range("E14").Value = "  -11.61%  "

$ws.Range("E15").Value = "  -1.89%  "

$ws.Range("D16").Value = "3.398.75"
$ws.Range("E16").Value = "  -3.49%  "

$ws.Range("D17").Value = "64.792.34"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.84"
$ws.Range("E18").Value = "  -8.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.42"
$ws.Range("E19").Value = "  -14.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.76"
$ws.Range("E20").Value = "  -6.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("E21").Value = "  -5.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.78"
$ws.Range("E22").Value = "  -8.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.546"
$ws.Range("E23").Value = "  -8.38%  "

$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.48"
$ws.Range("E25").Value = "  -7.63%  "

$ws.Range("D26").Value = "3.521.37"
$ws.Range("E26").Value = "  -3.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000102"
$ws.Range("E27").Value = "  -11.22%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -10.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  -10.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -10.80%  "

$ws.Range("D32").Value = "3.397.58"
$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -7.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.62"
$ws.Range("E35").Value = "  -6.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.85"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.56"
$ws.Range("E37").Value = "  -12.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.12"
$ws.Range("E38").Value = "  -12.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").Value = "  -7.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.61"
$ws.Range("E40").Value = "  -12.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0742"
$ws.Range("E41").Value = "  -8.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("E42").Value = "  -5.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.14"
$ws.Range("E43").Value = "  -4.68%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.28"
$ws.Range("E45").Value = "  -15.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.57"
$ws.Range("E46").Value = "  -10.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.08"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.67"
$ws.Range("E48").Value = "  -4.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.39"
$ws.Range("E49").Value = "  -8.49%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.138.86"
$ws.Range("E50").Value = "  -9.08%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("E51").Value = "  -15.86%  "
